# Growth-data edit: the Sept 6 density readings were originally split across
# two columns ("sept_6_density1" in F and "sept_6_densityF" in G). The
# author consolidated them into a single "sept_6_density" column, removing
# the old F column (sept_6_density1) entirely and keeping what used to be
# column G (sept_6_densityF) as the new column F, renamed "sept_6_density".
# Everything to the right (G:O, i.e. sept_6_quadrat .. sept_16_zj) shifts
# one column to the left (G:N) to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "sept_6_density1" column (F). This shifts the old
# "sept_6_densityF" column (G) into F, and everything after it (H:O)
# one step left (G:N) as well - exactly matching the consolidated layout.
$ws.Columns("F:F").Delete()

# Rename the (now merged) column header to the new combined field name.
$ws.Range("F1").Value = "sept_6_density"

# The column was resized to fit the new header text ("sept_6_density" is
# 14 characters). ColumnWidth uses Excel's "characters" unit, which is
# offset from the stored/serialized width by the default column padding
# (5/6 of a character), so we back that out to land on a stored width of 14.
$ws.Columns("F:F").ColumnWidth = 13.166666666666666

# Reflect the author's final selection/cursor position on the sheet.
$ws.Range("F1").Select()
